$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3: replace Gagana Puli with Aditya Jandhyala, role -> Admin ---
$ws.Range("B3").Value = "aditya.jandhyala@epsoftinc.com"
$ws.Range("D3").Value = "Admin"

# --- Add new row 4: Salma Shaik ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "salma.shaik@epsoftinc.com"
$ws.Range("C4").Value = "2.0"
$ws.Range("D4").Value = "Admin"

# --- Rebuild hyperlinks (the runtime's Hyperlinks.Delete() clears the whole
#     sheet collection, so remove them all and re-add the three needed) ---
[void]$ws.Range("A1").Hyperlinks.Delete()

[void]$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:lakshmi.kadali@epsoftinc.com", "", "mailto:lakshmi.kadali@epsoftinc.com")
[void]$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:aditya.jandhyala@epsoftinc.com")
[void]$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:salma.shaik@epsoftinc.com")

# Adding hyperlinks resets direct formatting, so restore the original look:
# B2 keeps the purple "visited" hyperlink style, B3/B4 use the plain blue one.
$ws.Range("B2").Font.Color = 8388736
$ws.Range("B2").Font.Underline = 2
$ws.Range("B2").HorizontalAlignment = -4108

$ws.Range("B3").Font.Color = 16711680
$ws.Range("B3").Font.Underline = 2
$ws.Range("B3").HorizontalAlignment = -4108

$ws.Range("B4").Font.Color = 16711680
$ws.Range("B4").Font.Underline = 2
$ws.Range("B4").HorizontalAlignment = -4108

# --- Match the author's final cursor position ---
[void]$ws.Range("D11").Select()
